$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values per row (columns D,H,I,J are numbers; E,F,G,K,L are text-formatted numbers)
$data = @{
    2  = @{ D=170; E="105,953,332.00"; F="496,525,037.00"; G="21.34"; H=27.298999999999999; I=700; J=103; K="14.66"; L="14.71" }
    3  = @{ D=170; E="101,914,515.00"; F="530,700,183.00"; G="19.20"; H=25.963999999999999; I=732; J=103; K="12.79"; L="14.07" }
    4  = @{ D=196; E="120,951,849.00"; F="516,061,290.00"; G="23.44"; H=25.873000000000001; I=728; J=110; K="14.42"; L="15.11" }
    5  = @{ D=156; E="95,699,464.00";  F="492,928,251.00"; G="19.41"; H=16.292000000000002; I=704; J=90;  K="13.87"; L="12.78" }
    6  = @{ D=155; E="100,426,431.00"; F="547,851,088.00"; G="18.33"; H=29.696000000000002; I=729; J=77;  K="11.35"; L="10.56" }
    7  = @{ D=170; E="136,052,680.00"; F="536,509,428.00"; G="25.36"; H=29.332000000000001; I=732; J=85;  K="13.47"; L="11.61" }
    8  = @{ D=181; E="122,646,112.00"; F="562,903,350.00"; G="21.79"; H=49.555;              I=732; J=108; K="14.13"; L="14.75" }
    9  = @{ D=139; E="95,842,398.00";  F="551,480,411.00"; G="17.38"; H=29.443000000000001; I=735; J=72;  K="10.63"; L="9.80"  }
    10 = @{ D=209; E="128,885,648.00"; F="534,643,664.00"; G="24.11"; H=37.000999999999998; I=731; J=112; K="13.98"; L="15.32" }
    11 = @{ D=185; E="126,921,858.00"; F="546,168,257.00"; G="23.24"; H=35.308999999999997; I=729; J=112; K="13.54"; L="15.36" }
    12 = @{ D=166; E="108,824,115.00"; F="541,797,311.00"; G="20.09"; H=32.646999999999998; I=732; J=98;  K="12.95"; L="13.39" }
    13 = @{ D=168; E="113,277,619.00"; F="529,573,730.00"; G="21.39"; H=33.149000000000001; I=736; J=98;  K="14.80"; L="13.32" }
    14 = @{ D=200; E="118,485,834.00"; F="529,964,396.00"; G="22.36"; H=40.856000000000002; I=738; J=93;  K="11.25"; L="12.60" }
    15 = @{ D=167; E="132,911,108.00"; F="535,710,779.00"; G="24.81"; H=54.57;              I=730; J=77;  K="13.84"; L="10.55" }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
    $ws.Range("H$row").Value = $vals.H
    $ws.Range("I$row").Value = $vals.I
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
}
